# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 16:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1730027
$ws.Range("C4").Value = 4752
$ws.Range("D4").Value = 480320
$ws.Range("E4").Value = 1148946
$ws.Range("G4").Value = 189
$ws.Range("H4").Value = 100761

# Row 70 - Azerbaiyan
$ws.Range("B70").Value = 4568
$ws.Range("C70").Value = 165
$ws.Range("D70").Value = 2897
$ws.Range("E70").Value = 1617
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 54

# Row 103 - Sri Lanka
$ws.Range("B103").Value = 1372
$ws.Range("C103").Value = 53
$ws.Range("E103").Value = 630

# Row 138 - Reunion
$ws.Range("B138").Value = 460
$ws.Range("C138").Value = 1
$ws.Range("E138").Value = 48
